$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 336, shifting existing rows 336:359 down to 337:360
$ws.Rows.Item(336).Insert()

# Fill in the new row 336 with the new weekly record
$r = 336
$ws.Cells.Item($r, 1).Value2 = 3
$ws.Cells.Item($r, 2).Value2 = "Femacal de La Calera"
$ws.Cells.Item($r, 3).Value2 = "Coquimbo"
$ws.Cells.Item($r, 4).Value2 = 44714
$ws.Cells.Item($r, 5).Value2 = 5
$ws.Cells.Item($r, 6).Value2 = 100112043
$ws.Cells.Item($r, 7).Value2 = "Pepino ensalada"
$ws.Cells.Item($r, 8).Value2 = "Sin especificar"
$ws.Cells.Item($r, 9).Value2 = "Primera"
$ws.Cells.Item($r, 10).Value2 = 130
$ws.Cells.Item($r, 11).Value2 = 19000
$ws.Cells.Item($r, 12).Value2 = 20000
$ws.Cells.Item($r, 13).Value2 = 19462
$ws.Cells.Item($r, 14).Value2 = "$/caja 70 unidades"
$ws.Cells.Item($r, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item($r, 16).Value2 = 278
$ws.Cells.Item($r, 17).Value2 = 70
$ws.Cells.Item($r, 18).Value2 = "Hortaliza"

# Ensure the date cell keeps the date/time number format used by column D elsewhere
$ws.Cells.Item($r, 4).NumberFormat = $ws.Cells.Item($r - 1, 4).NumberFormat

$wb.Save()
